# Fruta / hortaliza, semanal
# Insert the latest weekly price-report row for "Macroferia Regional de Talca"
# (Arándano (blue)) as the new row 141, pushing the existing rows 141:152
# down to 142:153.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 141; Excel shifts rows 141:152
# down to 142:153 and copies the formatting (incl. the date number format
# on column D) from the row above, same as a manual "Insert Row" in the UI.
$ws.Rows(141).Insert()

# Populate the new week's record.
$ws.Range("A141").Value = 5
$ws.Range("B141").Value = "Macroferia Regional de Talca"
$ws.Range("C141").Value = "Maule"
$ws.Range("D141").Value = 45265
$ws.Range("E141").Value = 7
$ws.Range("F141").Value = "Fruta"
$ws.Range("G141").Value = 100101
$ws.Range("H141").Value = "Berries"
$ws.Range("I141").Value = 100101001
$ws.Range("J141").Value = "Arándano (blue)"
$ws.Range("K141").Value = "Sin especificar"
$ws.Range("L141").Value = "Primera"
$ws.Range("M141").Value = 190
$ws.Range("N141").Value = 4000
$ws.Range("O141").Value = 4000
$ws.Range("P141").Value = 4000
$ws.Range("Q141").Value = "`$/bandeja 2 kilos"
$ws.Range("R141").Value = "Provincia de Curicó"
$ws.Range("S141").Value = 2000
$ws.Range("T141").Value = 2
